# SWT301_TeamEvaluation.xlsx update
# - Rename sheet "Assignment" -> "Lab2"
# - Update existing rows' Task text and the 2nd row's date
# - Append two new student rows (SE184409 / SE184761) with matching styling
# - Reflects the revised team evaluation data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Lab2"

# Row 4 (new): Lam Thi Ngoc Han
$ws.Range("B4").Value = "Lâm Thị Ngọc Hân"

# Row 5 (new): Pham Thi Thu Phuong
$ws.Range("B5").Value = "Phạm Thị Thu Phương"

$ws.Range("A4").Value = "SE184409"
$ws.Range("A5").Value = "SE184761"

# Task descriptions: new rows first share the "Confirm test case..." text,
# then the existing row 2 task is revised last
$ws.Range("C3").Value = "Confirm test case with web production & csv file, note bugs & test cases"
$ws.Range("C4").Value = "Confirm test case with web production & csv file, note bugs & test cases"
$ws.Range("C5").Value = "Confirm test case with web production & csv file, note bugs & test cases"
$ws.Range("C2").Value = "Builder Pattern implementation, Swimming Calorie Calculator Module, test case function & csv data implementation"

# Dates
$ws.Range("D3").Value = 45939
$ws.Range("D4").Value = 45939
$ws.Range("D5").Value = 45939

# Apply thin black borders to the new rows (matches the existing table styling)
$newRows = $ws.Range("A4:C5")
$newRows.Borders.Color = 0
$newRows.Borders.LineStyle = 1

# Date cells: copy the existing date-column format (border + number format) from D3
$ws.Range("D3").Copy()
$ws.Range("D4:D5").PasteSpecial(-4122)  # xlPasteFormats

# Widen FullName/Task columns to fit the longer revised content
$ws.Columns.Item(2).ColumnWidth = 22.0
$ws.Columns.Item(3).ColumnWidth = 98.5
